$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.482.53"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.571.59"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.76"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3706"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3376"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.151"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07542"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.18"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.024"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.970"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "1.571.28"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.48"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06772"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.353"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.43"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.24"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").Value = "22.482.23"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.380"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.634"
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.06"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.25"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.082"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.11"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "1.750.16"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.070"
$ws.Range("E32").Value = "  +8.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.205"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.014"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.804"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08360"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02482"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.358"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2302"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06548"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.438"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.30"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6228"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.13"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.805"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5864"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.08"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.075"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.239"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07312"
$ws.Range("E51").Value = "  -0.23%  "
